$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("One To One")

# Replace the old script function name with the new summary function name
# across every cell in column C that references it.
$null = $ws.Cells.Replace("forest_calcs.create_general_description_level", "01_general_description_summary - plot & level", 1)

# Switch the AutoFilter on column A (Script Group) from "General Descriptive" to "Health".
# This also updates which rows are hidden/visible to match the new filter criteria.
$null = $ws.Range("A1:L223").AutoFilter(1, @("Health"), 7)

# Update the active selection to reflect where the user ended up after filtering.
$null = $ws.Range("E240").Select()
